$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (order changed + renames)
$ws.Range("A1").Value = "Service Name"
$ws.Range("B1").Value = "Asset Group Name"
$ws.Range("C1").Value = "Asset  Name"
$ws.Range("D1").Value = "Asset Component Name"
$ws.Range("E1").Value = "Asset Owner Dept."
$ws.Range("F1").Value = "Asset Physical Location"
$ws.Range("G1").Value = "Asset Logical Location"
